$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 140, shifting existing rows 140:243 down to 141:244
$ws.Rows("140:140").Insert()

# Populate the newly inserted row 140 with the new record's data.
# Columns A,B,C,E,F,G,H,I,J,R are constant across this dataset (same market/region/product template).
$ws.Range("A140").Value = 5
$ws.Range("B140").Value = "Macroferia Regional de Talca"
$ws.Range("C140").Value = "Maule"
$ws.Range("D140").Value = 44673
$ws.Range("E140").Value = 7
$ws.Range("F140").Value = "Fruta"
$ws.Range("G140").Value = 100108
$ws.Range("H140").Value = "Tropicales y subtropicales"
$ws.Range("I140").Value = 100108005
$ws.Range("J140").Value = "Piña"
$ws.Range("K140").Value = "Sin especificar"
$ws.Range("L140").Value = "Tercera"
$ws.Range("M140").Value = 220
$ws.Range("N140").Value = 14000
$ws.Range("O140").Value = 15000
$ws.Range("P140").Value = 14455
$ws.Range("Q140").Value = "$/caja 16 unidades"
$ws.Range("R140").Value = "Ecuador"
$ws.Range("S140").Value = 903
$ws.Range("T140").Value = 16
